$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to match repulled data
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -8
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -7
$ws.Range("F11").Value = -8
$ws.Range("F13").Value = -4
$ws.Range("F14").Value = 3
$ws.Range("F16").Value = 8
$ws.Range("F17").Value = -2
$ws.Range("F20").Value = 0
